$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 21423
$ws1.Range("F3").Value = 3346
$ws1.Range("F4").Value = 857
$ws1.Range("F5").Value = 624
$ws1.Range("F6").Value = 554
$ws1.Range("F7").Value = 802
$ws1.Range("F8").Value = 302
$ws1.Range("F11").Value = 142
$ws1.Range("F12").Value = 578
$ws1.Range("F14").Value = 361
$ws1.Range("F15").Value = 39
$ws1.Range("F16").Value = 462
$ws1.Range("F17").Value = 200
$ws1.Range("F18").Value = 44
$ws1.Range("F20").Value = 85
$ws1.Range("F21").Value = 157

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 150
$ws2.Range("F10").Value = 171

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6177
$ws3.Range("F3").Value = 733
$ws3.Range("F4").Value = 731
$ws3.Range("F5").Value = 1737
$ws3.Range("F6").Value = 90

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6177
$ws4.Range("F3").Value = 733
$ws4.Range("F4").Value = 731
$ws4.Range("F5").Value = 1737
$ws4.Range("F6").Value = 21423
$ws4.Range("F7").Value = 3346
$ws4.Range("F8").Value = 857
$ws4.Range("F9").Value = 150
$ws4.Range("F10").Value = 90
$ws4.Range("F11").Value = 624
$ws4.Range("F12").Value = 554
$ws4.Range("F13").Value = 802
$ws4.Range("F14").Value = 302
$ws4.Range("F20").Value = 142
$ws4.Range("F23").Value = 578
$ws4.Range("F27").Value = 361
$ws4.Range("F28").Value = 171
$ws4.Range("F29").Value = 39
$ws4.Range("F30").Value = 462
$ws4.Range("F32").Value = 200
$ws4.Range("F33").Value = 44
$ws4.Range("F37").Value = 85
$ws4.Range("F43").Value = 157
